# Update publications and patents:
# Rows 14 and 15 (the "authors" column, A) contained the same author list
# but typed with full-width Chinese commas ("，"). Normalize them to use
# regular half-width commas (",") to match the style used elsewhere in
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAuthors = "金日初,刘江,胡衍,缪函霈,姜泓羊,王星月,曾娜,叶海礼"

$ws.Range("A14").Value = $newAuthors
$ws.Range("A15").Value = $newAuthors

# Reflect the scrolled/selected state that was saved with the workbook.
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
$ws.Range("A15").Select()
